$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "hsBJT120"
$ws.Range("B2").Value = 23091323
$ws.Range("C2").Value = "rcuznfs75"
$ws.Range("D2").Value = "Qp4T&!n9"
$ws.Range("F2").Value = "kJCJKswv"
$ws.Range("G2").Value = "inor"
